# Fix booleans: clear the redundant "is_comment"/"exportable" cells that
# held the literal string "False" (the default), and drop the orphaned
# "False" entry from the shared-strings table by simply not referencing it
# anymore. Rows whose comment flag is True (rows 2 and 6) are left intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# is_comment column (C) was "False" on these data rows -> clear it
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()

# row 9 (menu.share.not.exported) was not exportable -> clear both
# exportable (B) and is_comment (C)
$ws.Range("B9").ClearContents()
$ws.Range("C9").ClearContents()
